# Excel controller now uses ParseExact instead of Convert.ToDateTime when
# formatting appointment timestamps, so the formatted strings no longer
# carry a ":00" seconds component. Update the three date/time text cells
# accordingly and move the active selection down from A5 to A4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "15/03/2019 5:55 PM"
$ws.Range("A2").Value = "22/03/2019 12:52 PM"
$ws.Range("A3").Value = "23/04/2019 12:42 PM"

$ws.Range("A4").Select() | Out-Null
